$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure the "TrainDetails" block (old cols A:G, rows 8-12) into two tables:
#     "TravelDetails" (TravelID/TrainNo/ClassId/QuotaId, cols C:F) and a new
#     "TrainDetails" (TrainNo/TrainName/Source/Destination/Date, cols Q:U).

# Copy the existing header cell format (fillId=2 shaded header look) from C8
# onto the new TrainDetails header cell (S8, which sits outside the old A:G
# range that is about to be cleared).
$ws.Range("C8").Copy($ws.Range("S8"))

# Clear the old table (content + formatting) now that its layout/header moved.
# Note E8 (the new TravelDetails header) is inside this old A:G range, so its
# shaded-header look is (re)applied afterward.
$ws.Range("A8:G12").Clear()

# Give the TravelDetails header (E8) the same shaded header look used
# elsewhere on the sheet (copy format only, from another header cell).
$ws.Range("L8").Copy($ws.Range("E8"))

# New header captions.
$ws.Range("E8").Value() = "TravelDetails"
$ws.Range("S8").Value() = "TrainDetails"

# The new TrainDetails header gets its own distinct fill color (not the
# shared header style used elsewhere on the sheet).
$ws.Range("S8").Interior.Color = 15189684

# --- TravelDetails table (cols C:F) ---
$ws.Range("C9").Value() = "TravelID"
$ws.Range("D9").Value() = "TrainNo"
$ws.Range("E9").Value() = "ClassId"
$ws.Range("F9").Value() = "QuotaId"

$ws.Range("C10").Value() = 1
$ws.Range("D10").Value() = 4567
$ws.Range("E10").Value() = 2
$ws.Range("F10").Value() = 3

$ws.Range("C11").Value() = 2
$ws.Range("D11").Value() = 8987
$ws.Range("E11").Value() = 1
$ws.Range("F11").Value() = 2

$ws.Range("C12").Value() = 3
$ws.Range("D12").Value() = 9897
$ws.Range("E12").Value() = 2
$ws.Range("F12").Value() = 2

# --- New TrainDetails table (cols Q:U) ---
$ws.Range("Q9").Value() = "TrainNo"
$ws.Range("R9").Value() = "TrainName"
$ws.Range("S9").Value() = "Source"
$ws.Range("T9").Value() = "Destination"
$ws.Range("U9").Value() = "Date"

$ws.Range("Q10").Value() = 4567
$ws.Range("R10").Value() = "Kovai express"
$ws.Range("S10").Value() = "chennai"
$ws.Range("T10").Value() = "coimbatore"
$ws.Range("U10").Value() = "15.04.22"

$ws.Range("Q11").Value() = 8987
$ws.Range("R11").Value() = "keralaexpress"
$ws.Range("S11").Value() = "chennai"
$ws.Range("T11").Value() = "kozhikode"
$ws.Range("U11").Value() = "13.4.22"

$ws.Range("Q12").Value() = 9897
$ws.Range("R12").Value() = "lucknowexpress"
$ws.Range("S12").Value() = "chennai"
$ws.Range("T12").Value() = "lucknow"
$ws.Range("U12").Value() = "13.4.22"

# The Payment table (cols J:M) is untouched, and the rest of the sheet
# (class/quota/preference/ticket tables below) is also untouched.

# Move the active selection, matching where editing ended up.
$ws.Range("Y42").Select() | Out-Null
